# The "References" slides (7, 13, 19, 25) each list four reference links in
# the content placeholder, one per paragraph (after a leading blank
# paragraph). The edit rotates that list left by one position:
#   en.wikipedia.org -> www.nih.gov -> scholar.google.com -> www.jstor.org -> (back to) en.wikipedia.org
# i.e. each URL is replaced by the URL that originally followed it, and the
# first URL wraps around to become the last.

$p = $ppt.ActivePresentation

$newUrls = @(
    "https://www.nih.gov/",
    "https://scholar.google.com/",
    "https://www.jstor.org/",
    "https://en.wikipedia.org/wiki/Main_Page"
)

$slideIndexes = @(7, 13, 19, 25)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(2)
    $tr = $shape.TextFrame.TextRange

    # Paragraph 1 is the blank line before the links; paragraphs 2-5 hold the
    # four URLs, in order.
    for ($i = 0; $i -lt $newUrls.Count; $i++) {
        $para = $tr.Paragraphs($i + 2)
        $para.Runs(1).Text = $newUrls[$i]
    }
}
